$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows starting at row 47, pushing the existing MapNpcData rows
# (currently rows 47-48) down to rows 50-51.
$ws.Range("A47:D49").Insert()

# New LoadingData rows (47-49)
$ws.Range("A47").Value = "LoadingData.Loading.1000027"
$ws.Range("B47").Value = "솜사탕상점 로딩화면"

$ws.Range("A48").Value = "LoadingData.Loading.1000028"
$ws.Range("B48").Value = "파우더상점 로딩화면"

$ws.Range("A49").Value = "LoadingData.Loading.1000029"
$ws.Range("B49").Value = "이불상점 로딩화면"
